$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.817.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.00%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.672.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.76%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '325.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.67%  '

# Row 7
$ws.Range("E7").Value = '  +1.87%  '

# Row 8
$ws.Range("E8").Value = '  +0.09%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.552'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.10%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0822'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.00%  '

# Row 13
$ws.Range("E13").Value = '  +0.27%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.76%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.092.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.73%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.661.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.19%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.874'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.805.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.23%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.48%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.96%  '

# Row 21
$ws.Range("E21").Value = '  +1.36%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0961'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.79%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '276.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.38%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.63%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.40%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.73%  '

# Row 27
$ws.Range("E27").Value = '  +0.04%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.42%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.33%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.141'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.65%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.71%  '

# Row 33
$ws.Range("E33").Value = '  +4.57%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0807'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.32%  '

# Row 36
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.01%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.04'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.18%  '

# Row 38
$ws.Range("E38").Value = '  +7.34%  '

# Row 39
$ws.Range("E39").Value = '  +9.58%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.53%  '

# Row 41
$ws.Range("E41").Value = '  +2.04%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.95%  '

# Row 43
$ws.Range("E43").Value = '  -0.03%  '

# Row 44
$ws.Range("E44").Value = '  +6.70%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.119.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.78%  '

# Row 46
$ws.Range("E46").Value = '  +6.67%  '

# Row 47
$ws.Range("E47").Value = '  +8.99%  '

# Row 48
$ws.Range("E48").Value = '  +7.39%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.58%  '

# Row 50
$ws.Range("E50").Value = '  +5.21%  '

# Row 51
$ws.Range("E51").Value = '  +7.01%  '
